$wb = $excel.ActiveWorkbook

# 1. Rename "Sheet2" -> "payment-request"
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "payment-request"

# 2. Fix the 'CASH' string (curly closing quote -> straight quote) in J3.
#    A leading apostrophe typed/assigned directly is interpreted by Excel as a
#    "force text" marker and stripped from the stored value (and flips on the
#    quote-prefix cell style), so instead we compute the literal text with a
#    formula in a scratch cell, paste its value back into J3, then clean up
#    the scratch cell. This keeps a literal leading apostrophe in the text
#    without touching J3's style.
$ws.Range("Z1").Formula = "=CHAR(39)&""CASH""&CHAR(39)"
$ws.Range("Z1").Copy()
$ws.Range("J3").PasteSpecial(-4163) # xlPasteValues
$ws.Range("Z1").ClearContents()

# 3. Add a new row (row 4) describing the column data types.
$ws.Range("A4").Value = "Long"
$ws.Range("B4").Value = "String"
$ws.Range("C4").Value = "Date"
$ws.Range("D4").Value = "Date"
$ws.Range("E4").Value = "String"
$ws.Range("F4").Value = "String"
$ws.Range("G4").Value = "String"
$ws.Range("H4").Value = "String"
$ws.Range("I4").Value = "Long"
$ws.Range("J4").Value = "String"
$ws.Range("K4").Value = "Float"

# 4. Update the selection/active cell on the "payment-request" sheet.
$ws.Activate()
$ws.Range("D28").Select()
